# Update sheet name to 'Chart'
#
# The workbook's only sheet is renamed from "Multiple Queries" to "Chart".
# The embedded chart (xl/charts/chart1.xml) has series formulas that refer
# to the worksheet by name (e.g. 'Multiple Queries'!$C$7); when a sheet is
# renamed in Excel, those references are automatically rewritten to use the
# new sheet name. Since "Chart" has no spaces it no longer needs the
# surrounding single quotes (Chart!$C$7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet.
$ws.Name = "Chart"

# Update the chart's plotted series so their source formulas point at the
# renamed sheet, same as Excel does automatically on a sheet rename.
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection(1)

$ser.Name = "=Chart!`$C`$7"
$ser.XValues = "=Chart!`$A`$8:`$A`$13"
$ser.Values = "=Chart!`$C`$8:`$C`$13"
